$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing header cells I1 and J1 (append _x suffix)
$ws.Range("I1").Value = "DS_ESTADO_CIVIL_x"
$ws.Range("J1").Value = "DS_GRAU_INSTRUCAO_x"

# Add new header cells K1:N1
$ws.Range("K1").Value = "DS_ESTADO_CIVIL_y"
$ws.Range("L1").Value = "DS_GRAU_INSTRUCAO_y"
$ws.Range("M1").Value = "SQ_CANDIDATO"
$ws.Range("N1").Value = "SG_PARTIDO"

# Apply the same style used by the other header cells (bold + border)
$ws.Range("J1").Copy()
$ws.Range("K1:N1").PasteSpecial(-4122)

# SG_PARTIDO values for each candidate row (2-84)
$parties = @{
    2 = "NOVO"
    3 = "NOVO"
    4 = "NOVO"
    5 = "REDE"
    6 = "PC do B"
    7 = "PC do B"
    8 = "PC do B"
    9 = "PC do B"
    10 = "PV"
    11 = "PV"
    12 = "PSTU"
    13 = "PROS"
    14 = "PP"
    15 = "PP"
    16 = "PP"
    17 = "PSC"
    18 = "PSC"
    19 = "PSB"
    20 = "PSB"
    21 = "PTC"
    22 = "PP"
    23 = "PRB"
    24 = "DEM"
    25 = "PSDB"
    26 = "PSDB"
    27 = "PSDB"
    28 = "PSL"
    29 = "PSL"
    30 = "PSL"
    31 = "PSL"
    32 = "PSL"
    33 = "PSL"
    34 = "PSL"
    35 = "PSL"
    36 = "PSL"
    37 = "PSL"
    38 = "PATRIOTA"
    39 = "PATRIOTA"
    40 = "PATRIOTA"
    41 = "PPL"
    42 = "PATRIOTA"
    43 = "PPL"
    44 = "PPL"
    45 = "PSL"
    46 = "PCO"
    47 = "PSOL"
    48 = "PSOL"
    49 = "PCB"
    50 = "PMN"
    51 = "AVANTE"
    52 = "PR"
    53 = "MDB"
    54 = "PT"
    55 = "PDT"
    56 = "PDT"
    57 = "PSD"
    58 = "PSD"
    59 = "PPS"
    60 = "PRTB"
    61 = "PRTB"
    62 = "PRTB"
    63 = "MDB"
    64 = "PTB"
    65 = "PTB"
    66 = "PODE"
    67 = "MDB"
    68 = "PHS"
    69 = "PHS"
    70 = "PHS"
    71 = "PHS"
    72 = "PHS"
    73 = "PHS"
    74 = "PHS"
    75 = "SOLIDARIEDADE"
    76 = "SOLIDARIEDADE"
    77 = "DC"
    78 = "DC"
    79 = "DEM"
    80 = "PTC"
    81 = "PP"
    82 = "PPL"
    83 = "PRTB"
    84 = "DEM"
}

for ($row = 2; $row -le 84; $row++) {
    $ws.Cells.Item($row, 11).Value = $ws.Cells.Item($row, 9).Text
    $ws.Cells.Item($row, 12).Value = $ws.Cells.Item($row, 10).Text
    $ws.Cells.Item($row, 13).Value = $ws.Cells.Item($row, 6).Value2
    $ws.Cells.Item($row, 14).Value = $parties[$row]
}
